$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:L2").Value = 33.33333333333333
$ws.Range("M2").Value = 0.1666666666666667
$ws.Range("N2").Value = 0.1111111111111111
$ws.Range("O2").Value = 0.3333333333333334
